$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.433.89"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'1.866.55"
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'0.7069"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.79%  "
$ws.Range("D6").Value = "'243.38"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "'0.3142"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.13%  "
$ws.Range("D9").Value = "'0.07861"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.23%  "
$ws.Range("D10").Value = "'24.49"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.49%  "
$ws.Range("D11").Value = "'0.08023"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.89%  "
$ws.Range("D12").Value = "'1.860.69"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.64%  "
$ws.Range("D13").Value = "'5.203"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.18%  "
$ws.Range("D15").Value = "'0.7011"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.39%  "
$ws.Range("D16").Value = "'6.455"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.43%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "'0.000008376"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.00%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "'29.435.70"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.20%  "
$ws.Range("D19").Value = "'252.70"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.86%  "
$ws.Range("D20").Value = "'2.129.96"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.89%  "
$ws.Range("E21").Value = "  -1.53%  "
$ws.Range("D23").Value = "'7.607"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.89%  "
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("D25").Value = "'0.1558"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.07%  "
$ws.Range("D26").Value = "'9.012"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.01%  "
$ws.Range("D27").Value = "'160.84"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.53%  "
$ws.Range("D28").Value = "'18.74"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.64%  "
$ws.Range("E29").Value = "  -0.57%  "
$ws.Range("D30").Value = "'4.326"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.54%  "
$ws.Range("D31").Value = "'4.286"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.55%  "
$ws.Range("D32").Value = "'1.211"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.47%  "
$ws.Range("D33").Value = "'0.05307"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.99%  "
$ws.Range("D34").Value = "'1.886"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.23%  "
$ws.Range("D35").Value = "'0.7526"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.68%  "
$ws.Range("D36").Value = "'1.167"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.67%  "
$ws.Range("D37").Value = "'2.711"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.71%  "
$ws.Range("D38").Value = "'0.01878"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.75%  "
$ws.Range("D39").Value = "'1.265.84"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.25%  "
$ws.Range("D40").Value = "'2.742"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.39%  "
$ws.Range("D41").Value = "'0.8981"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.44%  "
$ws.Range("D42").Value = "'109.27"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.07%  "
$ws.Range("D43").Value = "'5.965"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.36%  "
$ws.Range("D44").Value = "'71.38"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.40%  "
$ws.Range("E45").Value = "  -0.06%  "
$ws.Range("E46").Value = "  -1.24%  "
$ws.Range("D47").Value = "'2.029.59"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.41%  "
$ws.Range("D48").Value = "'1.790"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.10%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'9.547"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.12%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").Value = "'0.5175"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.92%  "
$ws.Range("D51").Value = "'0.4311"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.67%  "
